# Applies numeric corrections produced by the scheduled market-data runner
# to the currentAveragePrice / LevePrice / LeveProfit columns (H..N) across
# the 8 per-class leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 16667430
$ws.Range("I80").Value = 25000304
$ws.Range("K80").Value = 75000912
$ws.Range("M80").Value = -74999914
$ws.Range("H83").Value = 16667430
$ws.Range("I83").Value = 25000304
$ws.Range("K83").Value = 225002736
$ws.Range("M83").Value = -224997744
$ws.Range("H86").Value = 166669300
$ws.Range("I86").Value = 166669300
$ws.Range("K86").Value = 166669300
$ws.Range("M86").Value = -166668177
$ws.Range("H89").Value = 166669300
$ws.Range("I89").Value = 166669300
$ws.Range("K89").Value = 833346500
$ws.Range("M89").Value = -833340884
$ws.Range("H92").Value = 90909730
$ws.Range("I92").Value = 142857980
$ws.Range("K92").Value = 142857980
$ws.Range("M92").Value = -142856732
$ws.Range("H129").Value = 2780.3333
$ws.Range("J129").Value = 6298
$ws.Range("L129").Value = 18894
$ws.Range("N129").Value = -28894
$ws.Range("H132").Value = 4878.8623
$ws.Range("I132").Value = 1680.12
$ws.Range("K132").Value = 5040.36
$ws.Range("M132").Value = -2510.36
$ws.Range("H137").Value = 12719048
$ws.Range("J137").Value = 33337154
$ws.Range("L137").Value = 100011462
$ws.Range("N137").Value = -100016562
$ws.Range("H141").Value = 2497.5
$ws.Range("I141").Value = 2497.5
$ws.Range("K141").Value = 7492.5
$ws.Range("M141").Value = -2312.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20526.316
$ws.Range("I32").Value = 20173.264
$ws.Range("K32").Value = 20173.264
$ws.Range("M32").Value = -19886.264
$ws.Range("H61").Value = 5537.25
$ws.Range("I61").Value = 5956.857
$ws.Range("K61").Value = 5956.857
$ws.Range("M61").Value = -5744.857
$ws.Range("H74").Value = 19232304
$ws.Range("I74").Value = 20834830
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 20834830
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -20833956
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 19232304
$ws.Range("I77").Value = 20834830
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 104174150
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -104169782
$ws.Range("N77").Value = -18736
$ws.Range("H122").Value = 3716.5122
$ws.Range("I122").Value = 3629.1025
$ws.Range("K122").Value = 10887.3075
$ws.Range("M122").Value = -8437.307499999999
$ws.Range("H132").Value = 1887.9762
$ws.Range("I132").Value = 1805.0256
$ws.Range("K132").Value = 5415.0768
$ws.Range("M132").Value = -2885.0768
$ws.Range("H136").Value = 5537.25
$ws.Range("I136").Value = 5956.857
$ws.Range("K136").Value = 17870.571
$ws.Range("M136").Value = -15320.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 45375
$ws.Range("J40").Value = 45375
$ws.Range("L40").Value = 45375
$ws.Range("N40").Value = -45905
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H99").Value = 2419.25
$ws.Range("I99").Value = 2557.1667
$ws.Range("K99").Value = 2557.1667
$ws.Range("M99").Value = -1059.1667
$ws.Range("H105").Value = 4846.5
$ws.Range("I105").Value = 4688.2144
$ws.Range("K105").Value = 4688.2144
$ws.Range("M105").Value = -2941.2144
$ws.Range("H134").Value = 2730.9722
$ws.Range("I134").Value = 2014.5454
$ws.Range("K134").Value = 6043.6362
$ws.Range("M134").Value = -3508.6362
$ws.Range("H137").Value = 96666.336
$ws.Range("J137").Value = 97500
$ws.Range("L137").Value = 97500
$ws.Range("N137").Value = -107700
$ws.Range("H138").Value = 55554
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 109207.25
$ws.Range("J139").Value = 109554.14
$ws.Range("L139").Value = 109554.14
$ws.Range("N139").Value = -119834.14

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 815.73334
$ws.Range("I16").Value = 715.2857
$ws.Range("K16").Value = 715.2857
$ws.Range("M16").Value = -428.2857
$ws.Range("H107").Value = 907.5
$ws.Range("J107").Value = 1060.6666
$ws.Range("L107").Value = 1060.6666
$ws.Range("N107").Value = -4900.6666
$ws.Range("H113").Value = 815.73334
$ws.Range("I113").Value = 715.2857
$ws.Range("K113").Value = 715.2857
$ws.Range("M113").Value = 1454.7143
$ws.Range("H132").Value = 83341270
$ws.Range("I132").Value = 95239920
$ws.Range("J132").Value = 50750
$ws.Range("K132").Value = 285719760
$ws.Range("L132").Value = 152250
$ws.Range("M132").Value = -285717230
$ws.Range("N132").Value = -157310
$ws.Range("H141").Value = 119357.65
$ws.Range("I141").Value = 83600
$ws.Range("J141").Value = 123515.51
$ws.Range("K141").Value = 83600
$ws.Range("L141").Value = 123515.51
$ws.Range("M141").Value = -78420
$ws.Range("N141").Value = -133875.51

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 468.70587
$ws.Range("J107").Value = 560.6667
$ws.Range("L107").Value = 1682.0001
$ws.Range("N107").Value = -5522.0001
$ws.Range("H132").Value = 1115.8
$ws.Range("I132").Value = 909.90625
$ws.Range("K132").Value = 8189.15625
$ws.Range("M132").Value = -5659.15625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 36692
$ws.Range("J62").Value = 36692
$ws.Range("L62").Value = 36692
$ws.Range("N62").Value = -38064
$ws.Range("H65").Value = 36692
$ws.Range("J65").Value = 36692
$ws.Range("L65").Value = 110076
$ws.Range("N65").Value = -116940
$ws.Range("H68").Value = 20000
$ws.Range("J68").Value = 20000
$ws.Range("L68").Value = 20000
$ws.Range("N68").Value = -21622
$ws.Range("H71").Value = 20000
$ws.Range("J71").Value = 20000
$ws.Range("L71").Value = 60000
$ws.Range("N71").Value = -68112
$ws.Range("H80").Value = 120019.6
$ws.Range("I80").Value = 206399.4
$ws.Range("J80").Value = 33639.8
$ws.Range("K80").Value = 206399.4
$ws.Range("L80").Value = 33639.8
$ws.Range("M80").Value = -205401.4
$ws.Range("N80").Value = -35635.8
$ws.Range("H83").Value = 120019.6
$ws.Range("I83").Value = 206399.4
$ws.Range("J83").Value = 33639.8
$ws.Range("K83").Value = 1031997
$ws.Range("L83").Value = 168199
$ws.Range("M83").Value = -1027005
$ws.Range("N83").Value = -178183
$ws.Range("H122").Value = 225473.95
$ws.Range("I122").Value = 358557.94
$ws.Range("K122").Value = 1075673.82
$ws.Range("M122").Value = -1073223.82
$ws.Range("H132").Value = 113627.164
$ws.Range("I132").Value = 154999.16
$ws.Range("K132").Value = 464997.48
$ws.Range("M132").Value = -462467.48
$ws.Range("H141").Value = 33030.832
$ws.Range("J141").Value = 33030.832
$ws.Range("L141").Value = 33030.832
$ws.Range("N141").Value = -43390.832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4621.4707
$ws.Range("I46").Value = 2948
$ws.Range("J46").Value = 4844.6
$ws.Range("K46").Value = 2948
$ws.Range("L46").Value = 4844.6
$ws.Range("M46").Value = -2760
$ws.Range("N46").Value = -5220.6
$ws.Range("H93").Value = 1781.2667
$ws.Range("I93").Value = 1781.2667
$ws.Range("K93").Value = 1781.2667
$ws.Range("M93").Value = -533.2666999999999
$ws.Range("H122").Value = 7577
$ws.Range("I122").Value = 4338.9443
$ws.Range("K122").Value = 13016.8329
$ws.Range("M122").Value = -10566.8329
$ws.Range("H132").Value = 2942.2627
$ws.Range("I132").Value = 2865.897
$ws.Range("J132").Value = 3109.7742
$ws.Range("K132").Value = 8597.690999999999
$ws.Range("L132").Value = 9329.3226
$ws.Range("M132").Value = -6067.690999999999
$ws.Range("N132").Value = -14389.3226
$ws.Range("H136").Value = 4260.375
$ws.Range("I136").Value = 2828.6453
$ws.Range("J136").Value = 6035.72
$ws.Range("K136").Value = 8485.9359
$ws.Range("L136").Value = 18107.16
$ws.Range("M136").Value = -5935.9359
$ws.Range("N136").Value = -23207.16

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2471.48
$ws.Range("I122").Value = 868.6923
$ws.Range("K122").Value = 2606.0769
$ws.Range("M122").Value = -156.0769
$ws.Range("H132").Value = 66669824
$ws.Range("I132").Value = 83335890
$ws.Range("J132").Value = 5555
$ws.Range("K132").Value = 250007670
$ws.Range("L132").Value = 16665
$ws.Range("M132").Value = -250005140
$ws.Range("N132").Value = -21725
$ws.Range("H136").Value = 3632.027
$ws.Range("I136").Value = 2585.5356
$ws.Range("K136").Value = 7756.6068
$ws.Range("M136").Value = -5206.6068
